# Apply updated cryptocurrency price/volume data per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.553.87"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "2.102.28"
$ws.Range("E3").Value = "  +9.76%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.51"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("E6").Value = "  -6.53%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.89"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.04%  "
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("E11").Value = "  -2.40%  "
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.62"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").Value = "2.407.19"
$ws.Range("E14").Value = "  +9.69%  "
$ws.Range("E15").Value = "  +2.60%  "
$ws.Range("D16").Value = "2.097.00"
$ws.Range("E16").Value = "  +9.40%  "
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").Value = "36.585.70"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.86"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("D20").Value = "0.0₃0832"
$ws.Range("E20").Value = "  -3.27%  "
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "240.71"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.23%  "
$ws.Range("E23").Value = "  +1.25%  "
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("E25").Value = "  -4.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.42"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.26"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +13.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.16"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -9.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.82"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +58.13%  "
$ws.Range("E31").Value = "  -5.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.49"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.63%  "
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("E34").Value = "  +22.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.986"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +12.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0901"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("E39").Value = "  -5.67%  "
$ws.Range("E40").Value = "  -11.40%  "
$ws.Range("E41").Value = "  +6.57%  "
$ws.Range("E42").Value = "  -1.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "98.43"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.53%  "
$ws.Range("E44").Value = "  -5.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.06"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -9.16%  "
$ws.Range("D46").Value = "1.335.30"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0844"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.08"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +9.50%  "
$ws.Range("D49").Value = "2.307.97"
$ws.Range("E49").Value = "  +10.22%  "
$ws.Range("E50").Value = "  +1.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.24"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -6.00%  "
